$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date line
Replace-Text "2024-12-05 Thursday" "2024-12-06 Friday"

# Row 1
Replace-Text "33÷6=" "41÷8="
Replace-Text "23÷3=" "93÷5="
Replace-Text "73÷6=" "54÷4="
Replace-Text "43÷2=" "70÷7="
Replace-Text "80÷3=" "27÷8="

# Row 2 (note: process the "92÷8=" -> "65÷8=" cell before "77÷6=" -> "92÷8=" to
# avoid the new value of the first cell being re-matched by the third cell's rule)
Replace-Text "92÷8=" "65÷8="
Replace-Text "77÷6=" "92÷8="
Replace-Text "67÷6=" "31÷3="
Replace-Text "13÷2=" "16÷6="
Replace-Text "47÷3=" "28÷8="

# Row 3
Replace-Text "81÷3=" "24÷3="
Replace-Text "55÷7=" "63÷2="
Replace-Text "55÷4=" "23÷8="
Replace-Text "41÷9=" "94÷3="
Replace-Text "13÷6=" "21÷5="

# Row 4
Replace-Text "82÷6=" "87÷3="
Replace-Text "65÷5=" "38÷8="
Replace-Text "49÷6=" "43÷4="
Replace-Text "29÷3=" "97÷9="
Replace-Text "32÷2=" "66÷4="

# Row 5
Replace-Text "83÷9=" "24÷8="
Replace-Text "68÷3=" "16÷7="
Replace-Text "60÷6=" "27÷5="
Replace-Text "35÷2=" "45÷9="
Replace-Text "62÷8=" "40÷6="
